$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the "Haba" sheet. It belongs right
# above the (former) last data row, so insert a fresh row at 62 - this
# shifts the existing rows 62:83 down to 63:84 (dimension grows to R84) -
# and then populate the new row with its data.
$ws.Rows.Item(62).Insert()

$ws.Range("A62").Value = 10
$ws.Range("B62").Value = "Vega Modelo de Temuco"
$ws.Range("C62").Value = "La Araucanía"
$ws.Range("D62").Value = 44875
$ws.Range("E62").Value = 9
$ws.Range("F62").Value = 100112026
$ws.Range("G62").Value = "Haba"
$ws.Range("H62").Value = "Sin especificar"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 110
$ws.Range("K62").Value = 10000
$ws.Range("L62").Value = 10000
$ws.Range("M62").Value = 10000
$ws.Range("N62").Value = "$/saco 25 kilos"
$ws.Range("O62").Value = "Región Metropolitana"
$ws.Range("P62").Value = 400
$ws.Range("Q62").Value = 25
$ws.Range("R62").Value = "Hortaliza"
